$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2024-04-15 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-16 Tuesday", 2) | Out-Null

# Update the table cells (row-major, 20 rows x 5 cols)
$t = $d.Tables(1)
$values = @(
    "18+63=",
    "66-30=",
    "42+3=",
    "66+5=",
    "10+54=",
    "47-13=",
    "66-11=",
    "88-3=",
    "79-30=",
    "19-5=",
    "99-67=",
    "88-34=",
    "70-35=",
    "93+5=",
    "71+11=",
    "76-61=",
    "17+43=",
    "42+53=",
    "31+30=",
    "61-35=",
    "39-23=",
    "56+34=",
    "50+3=",
    "78+2=",
    "7+66=",
    "28+66=",
    "15-13=",
    "58-13=",
    "95-24=",
    "97+1=",
    "75-63=",
    "54+3=",
    "69-45=",
    "85-56=",
    "27+24=",
    "81-0=",
    "7+39=",
    "37+16=",
    "53+39=",
    "49-23=",
    "53-20=",
    "15+29=",
    "17-4=",
    "32-0=",
    "29-27=",
    "48-0=",
    "99-37=",
    "41-12=",
    "11+46=",
    "42+10=",
    "24+69=",
    "93-24=",
    "96-17=",
    "34-27=",
    "12+22=",
    "95-28=",
    "70-53=",
    "84-28=",
    "62-27=",
    "14+33=",
    "97-12=",
    "3+26=",
    "14+46=",
    "84+8=",
    "56-34=",
    "3+86=",
    "60+0=",
    "2+22=",
    "71+26=",
    "37-6=",
    "46+6=",
    "26+32=",
    "5+71=",
    "81-30=",
    "28-9=",
    "60-6=",
    "77-40=",
    "42-32=",
    "66-15=",
    "76-8=",
    "50+37=",
    "94-32=",
    "35+43=",
    "76-4=",
    "81-10=",
    "48-23=",
    "58-15=",
    "40-31=",
    "92-36=",
    "65+20=",
    "61-34=",
    "65-51=",
    "93-15=",
    "11+48=",
    "32+63=",
    "34-30=",
    "63-25=",
    "51+41=",
    "42+50=",
    "35+0="
)

$idx = 0
for ($row = 1; $row -le 20; $row++) {
    for ($col = 1; $col -le 5; $col++) {
        $t.Cell($row, $col).Range.Text = $values[$idx]
        $idx++
    }
}

Write-Output "Done: updated $idx cells"